$d = $word.ActiveDocument

# Add the new "Notes" paragraph style (custom style, based on Normal,
# quick-style, 10pt / half-points 20 font size) used for notes added to
# figures and tables.
$notes = $d.Styles.Add("Notes", 1)
$notes.BaseStyle = "Normal"
$notes.Font.Size = 10
$notes.QuickStyle = $true

# Switch the template placeholder paragraph from the "Code" style to the
# new "Notes" style.
$p = $d.Paragraphs.Item(1)
$p.Style = "Notes"
